# Applies the "proyecto terminado para presentacion" edit to the
# plantilla_importacio workbook: updates the TIPO DOC VENTA / CATEGORIAS /
# TIPO DOCUMENTO lists, tweaks some label casing, and clears the lingering
# cell selection left in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ImportacionGeneral")

# --- Column A: --- CATEGORIAS --- ---------------------------------------
$ws.Range("A3").Value = "platos"
$ws.Range("A5").Value = "entrada"
$ws.Range("A6").Value = "postre"
$ws.Range("A7").Value = "bebidas"

# --- Column H/I: --- TIPO DOCUMENTO --- ---------------------------------
$ws.Range("H3").Value = "DNI"
$ws.Range("H4").Value = "RUC"
$ws.Range("H5").Value = "carnet extrajero"
$ws.Range("I4").Value = 12
$ws.Range("I5").Value = 12

# --- Column M/N: --- TIPO DOC VENTA --- ---------------------------------
$ws.Range("M3").Value = "Boleta"
$ws.Range("N3").Value = "B"
$ws.Range("M4").Value = "Factura"
$ws.Range("N4").Value = "F"
$ws.Range("M5").Value = "nota de Venta"
$ws.Range("N5").Value = "N"

# Give the new column N content (B/F/N, header "Ciglas") a best-fit-style
# width like the other lookup columns.
$ws.Columns.Item(14).ColumnWidth = 4.721354166666667

# The saved selection (D7:D11) is no longer part of the sheet view once the
# workbook is resaved, so reset the active cell back to A1.
$ws.Range("A1").Select() | Out-Null

$wb.Save()
